$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A173").Value = 'Job Title: Tech Lead - Golang / AWS (Serverless Architecture)'
$ws.Range("B173").Value = 'https://www.dice.com/job-detail/988b81d8-4af0-40bc-b81b-d5502653280f?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang'
$ws.Range("C173").Value = 'McLean, Virginia'
$ws.Range("D173").Value = 'Contract, Third Party'
$ws.Range("E173").Value = '$70 - $80'
$ws.Range("F173").Value = 'Sagarsoft'

$ws.Range("A174").Value = 'Cloud Automation Engineer (AWS | Python/Go)'
$ws.Range("B174").Value = 'https://www.dice.com/job-detail/e29f8561-e75f-4fe4-b56f-92316d8d4b14?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang'
$ws.Range("C174").Value = 'Hybrid in Hartford, Connecticut'
$ws.Range("D174").Value = 'Contract'
$ws.Range("E174").Value = 'Depends on Experience'
$ws.Range("F174").Value = 'TechTalent Solutions LLC'
